# TC003_LoginDDT: add two more dynamic/unique login test rows with their
# own hyperlinked email addresses (mirrors the existing bhaskar@gmail.com
# row already on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new valid-login test data -------------------------------------
$ws.Range("A6").Value = "bhaskar4@gmail.com"
$ws.Range("B6").Value = "Test123"
$ws.Range("C6").Value = "Valid"

# --- Row 7: new valid-login test data -------------------------------------
$ws.Range("A7").Value = "bhaskar3@gmail.com"
$ws.Range("B7").Value = "Test12345"
$ws.Range("C7").Value = "Valid"

# Turn the two new e-mail addresses into mailto hyperlinks, same as the
# existing bhaskar@gmail.com entry in A2.
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:bhaskar4@gmail.com", "", "mailto:bhaskar4@gmail.com", "bhaskar4@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:bhaskar3@gmail.com", "", "", "bhaskar3@gmail.com") | Out-Null

# Leave the selection where the user finished entering the new data.
$ws.Range("A7:C7").Select() | Out-Null
